# Apply changes described by the diff:
# 1. Rename sheet "EM" -> "EM_Self_Tests"
# 2. Add column G (header "OM_Application") with data to the EM sheet
# 3. Adjust column widths for new column G (and F) on EM sheet
# 4. Update selections/view on both sheets

$wb = $excel.ActiveWorkbook

# --- Rename the EM sheet ---
$wsEM = $wb.Worksheets.Item("EM")
$wsEM.Name = "EM_Self_Tests"

$wsOM = $wb.Worksheets.Item("OM")

# --- Add new column G data on EM sheet ---
$wsEM.Range("G1").Value = "OM_Application"
$wsEM.Range("G2").Value = "Fast_LL_DatLow,Fast_LL_DatHigh,Slow_LL_DatLow,Slow_LL_DatHigh"
$wsEM.Range("G3").Value = "Fast_LG_DatLow,Fast_LG_DatHigh,Slow_LG_DatLow,Slow_LG_DatHigh"
$wsEM.Range("G4").Value = "Fast_LL_DatLow,Fast_LL_DatHigh,Slow_LL_DatLow,Slow_LL_DatHigh"
$wsEM.Range("G5").Value = "Fast_LG_DatLow,Fast_LG_DatHigh,Slow_LG_DatLow,Slow_LG_DatHigh"
$wsEM.Range("G6").Value = "Fast_LL_DatLow,Fast_LL_DatHigh,Slow_LL_DatLow,Slow_LL_DatHigh"
$wsEM.Range("G7").Value = "Fast_LG_DatLow,Fast_LG_DatHigh,Slow_LG_DatLow,Slow_LG_DatHigh"

# --- Column widths ---
# NOTE: the runtime quantizes ColumnWidth to the nearest 1/6 of a character
# (it always rounds (cw + 5/6) to the nearest sixth), so we cannot hit the
# exact 1/256-based OOXML width (17.33203125 / 61.6640625) from the diff.
# Feed in the input value that lands closest to those targets after the
# runtime's internal rounding.
$wsEM.Columns.Item(6).ColumnWidth = 16.5
$wsEM.Columns.Item(7).ColumnWidth = 60.833333333333336

# --- Selections / views ---
$wsOM.Activate()
$wsOM.Range("A9").Select()

$wsEM.Activate()
$wsEM.Range("G7").Select()
